$wb = $excel.ActiveWorkbook

# --- Sheet "Clientes": field/validation-message edits -----------------------
$clientes = $wb.Worksheets.Item("Clientes")

# Row 4 "Nombre": use %s placeholder instead of hard-coded length in the
# "too short" message.
$clientes.Range("K4").Value = "Muy corto, minimo %s"

# Row 5 "Correo": field is no longer mandatory (server-side bug fix for
# optional fields) and its min-length message now uses %s.
$clientes.Range("E5").Value = "no"
$clientes.Range("K5").Value = "Al menos %s"

# Row 6 "Tipo cliente" (radio): parsley-compatible constraint syntax.
$clientes.Range("L6").Value = "N==Natural||J==Juridica"

# Row 9 "Observaciones": split textarea ("Area") from regular text boxes,
# dropping the old "mucho texto" data type.
$clientes.Range("A9").Value = "Area"
$clientes.Range("C9").Value = "texto"

# Row 13 "Recibir publicidad" (checkbox): parsley-compatible constraint syntax.
$clientes.Range("L13").Value = "1==Si"

# New row 14: "Numero empleados" field.
$clientes.Range("A14").Value = "Caja"
$clientes.Range("B14").Value = "Numero empleados"
$clientes.Range("C14").Value = "numero"
$clientes.Range("F14").Value = 10
$clientes.Range("K14").Value = "Se acepta minimo %s empleados"

# --- Sheet view / selection updates -----------------------------------------
$config = $wb.Worksheets.Item("Config")
$config.Range("B5").Select()

$paises = $wb.Worksheets.Item("t_Paises")
$paises.Range("E5").Select()

$clientes.Range("B18").Select()
$clientes.Activate()
